$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update report generated timestamp
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:00 AM"

# Update total billed amount
$ws.Range("C8").Value = 1067.46

# Clear Scope ID # value
$ws.Range("G10").Value = ""

# Update per-line pricing values
$ws.Range("H16").Value = 632.4
$ws.Range("H17").Value = 435.06
$ws.Range("H18").Value = 1067.46
